$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.636.27"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.302.32"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("D8").Value = "3.303.31"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -4.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.64"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -12.66%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.93"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "3.843.04"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "3.307.17"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("D18").Value = "63.433.14"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.67"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "619.23"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -8.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.08"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.45"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.40%  "
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("D40").Value = "0.0₃0720"
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.876.25"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("E48").Value = "  -6.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.124"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.25%  "
